$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 107 (date 43431) - add start/end times
$ws.Range("B107").Value = 0.76458333333333339
$ws.Range("C107").Value = 0.99930555555555556

# Row 108 (date 43432) - add end time only (start stays 0)
$ws.Range("C108").Value = 0.40486111111111112

# Row 109 (date 43433) - add start/end times
$ws.Range("B109").Value = 0.4909722222222222
$ws.Range("C109").Value = 0.72083333333333333

# Row 110 (date 43434) - add start/end times
$ws.Range("B110").Value = 0.72430555555555554
$ws.Range("C110").Value = 0.99930555555555556

# Row 111 (date 43435) - correct start time, add end time
$ws.Range("B111").Value = 0
$ws.Range("C111").Value = 0.32013888888888892

# Row 112 (date 43436) - add new start/end times
$ws.Range("B112").Value = 0.41250000000000003
$ws.Range("C112").Value = 0.5229166666666667

# Row 113 (date 43437) - add new (zero) start/end times
$ws.Range("B113").Value = 0
$ws.Range("C113").Value = 0

# Row 114 (date 43438) - add new (zero) start/end times
$ws.Range("B114").Value = 0
$ws.Range("C114").Value = 0

# Update the active selection to match the edited area
$ws.Range("B115").Select() | Out-Null
